# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates the "K" column (column G) values for rows 2-36 on the active
# worksheet, replacing the old "Strike#" derived values with the newly
# calculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newKValues = @{
    2  = 5
    3  = 7
    4  = 1
    5  = 3
    6  = 4
    7  = 2
    8  = 3
    9  = 6
    10 = 9
    11 = 1
    12 = 1
    13 = 6
    14 = 7
    15 = 7
    16 = 13
    17 = 6
    18 = 3
    19 = 1
    20 = 1
    21 = 8
    22 = 5
    23 = 13
    24 = 1
    25 = 6
    26 = 4
    27 = 4
    28 = 0
    29 = 5
    30 = 4
    31 = 6
    32 = 5
    33 = 0
    34 = 3
    35 = 2
    36 = 1
}

foreach ($row in $newKValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newKValues[$row]
}
